$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Name und Vorname..." header text in A1 -------------------
# Content change: "Ausbildungsnachweis Nr.#idNr" -> "Ausbildungsnachweis Nr. #idNr"
# (a space was inserted) and the blank line before "Betrieblicher..." was
# removed. Rewriting this cell's value also has the side effect of moving
# this shared string to the end of the shared-string table (it becomes a
# brand new unique string, while the two untouched strings used by A5/C5
# shift down to fill the gap) - which is exactly the reordering the diff
# shows for sharedStrings.xml and for the A1/A5/C5 cell references.
$newHeaderText = "Name und Vorname des Auszubildenen #idName" + [char]10 + `
    "#idYear. Ausbildungsjahr" + [char]10 + `
    "Ausbildungsnachweis Nr. #idNr" + [char]10 + `
    "Für die Woche vom #idFirstDate bis #idLastDate." + [char]10 + `
    "Betrieblicher Funktionsberreich: #idDepartment"

$ws.Range("A1").Value = $newHeaderText

# --- Update the selected range on the sheet view ---------------------------
$ws.Range("A1:F1").Select()
